$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some column-D prices are plain decimal numbers (e.g. "1.002"); Excel would
# auto-convert a bare .Value assignment of such a string into a Number cell,
# while the source file stores them as literal text (inlineStr). Force those
# specific cells to Text format first so the assignment below keeps them as text.
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '27.977.75'
$ws.Range('E2').Value = '  -0.98%  '
$ws.Range('D3').Value = '1.742.32'
$ws.Range('E3').Value = '  -3.75%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').Value = '335.69'
$ws.Range('E5').Value = '  -1.02%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.20%  '
$ws.Range('D7').Value = '0.3758'
$ws.Range('E7').Value = '  -4.16%  '
$ws.Range('D8').Value = '0.3351'
$ws.Range('E8').Value = '  -4.29%  '
$ws.Range('D9').Value = '45.30'
$ws.Range('E9').Value = '  -6.83%  '
$ws.Range('D10').Value = '1.114'
$ws.Range('E10').Value = '  -5.80%  '
$ws.Range('D11').Value = '0.07200'
$ws.Range('E11').Value = '  -4.81%  '
$ws.Range('D12').Value = '1.003'
$ws.Range('E12').Value = '  +0.25%  '
$ws.Range('D13').Value = '22.36'
$ws.Range('E13').Value = '  +0.92%  '
$ws.Range('D14').Value = '6.155'
$ws.Range('E14').Value = '  -5.97%  '
$ws.Range('E15').Value = '  -1.11%  '
$ws.Range('D16').Value = '1.746.59'
$ws.Range('E16').Value = '  -3.73%  '
$ws.Range('E17').Value = '  -4.57%  '
$ws.Range('D18').Value = '0.06566'
$ws.Range('E18').Value = '  -2.22%  '
$ws.Range('D19').Value = '79.34'
$ws.Range('E19').Value = '  -7.07%  '
$ws.Range('E20').Value = '  +0.21%  '
$ws.Range('D21').Value = '16.82'
$ws.Range('E21').Value = '  -5.60%  '
$ws.Range('D22').Value = '6.233'
$ws.Range('E22').Value = '  -5.51%  '
$ws.Range('D23').Value = '27.982.35'
$ws.Range('E23').Value = '  -0.97%  '
$ws.Range('D24').Value = '11.64'
$ws.Range('E24').Value = '  -6.83%  '
$ws.Range('E25').Value = '  -0.20%  '
$ws.Range('D26').Value = '153.82'
$ws.Range('E26').Value = '  -0.74%  '
$ws.Range('D27').Value = '19.80'
$ws.Range('E27').Value = '  -7.64%  '
$ws.Range('D28').Value = '2.312'
$ws.Range('E28').Value = '  -8.72%  '
$ws.Range('D29').Value = '1.947.13'
$ws.Range('E29').Value = '  -3.60%  '
$ws.Range('D30').Value = '131.34'
$ws.Range('E30').Value = '  -3.74%  '
$ws.Range('D31').Value = '1.245'
$ws.Range('E31').Value = '  -16.20%  '
$ws.Range('D32').Value = '4.022'
$ws.Range('E32').Value = '  +0.00%  '
$ws.Range('D33').Value = '5.768'
$ws.Range('E33').Value = '  -10.40%  '
$ws.Range('D34').Value = '0.08713'
$ws.Range('E34').Value = '  -1.69%  '
$ws.Range('D35').Value = '12.16'
$ws.Range('E35').Value = '  -8.17%  '
$ws.Range('D36').Value = '0.6683'
$ws.Range('E36').Value = '  -3.80%  '
$ws.Range('E37').Value = '  -6.22%  '
$ws.Range('D38').Value = '0.06198'
$ws.Range('E38').Value = '  -5.67%  '
$ws.Range('D39').Value = '5.142'
$ws.Range('E39').Value = '  -6.61%  '
$ws.Range('E40').Value = '  -5.75%  '
$ws.Range('D41').Value = '1.212'
$ws.Range('E41').Value = '  -4.42%  '
$ws.Range('D42').Value = '1.436'
$ws.Range('E42').Value = '  -11.27%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '7.958'
$ws.Range('E43').Value = '  -7.31%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  +0.19%  '
$ws.Range('D45').Value = '13.79'
$ws.Range('E45').Value = '  -6.16%  '
$ws.Range('D46').Value = '3.823'
$ws.Range('E46').Value = '  -1.38%  '
$ws.Range('D47').Value = '0.6036'
$ws.Range('E47').Value = '  -6.33%  '
$ws.Range('D48').Value = '127.51'
$ws.Range('E48').Value = '  -3.40%  '
$ws.Range('D49').Value = '2.016'
$ws.Range('E49').Value = '  -7.15%  '
$ws.Range('D50').Value = '1.177'
$ws.Range('E50').Value = '  +0.39%  '
$ws.Range('D51').Value = '0.07113'
$ws.Range('E51').Value = '  -1.88%  '
